$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = 0.3
$wsSummary.Range("B6").Value = 29
$wsSummary.Range("B9").Value = 44.83

# --- Sheet: Strategy Status ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 29
$wsStatus.Range("G4").Value = 44.83

# --- Sheet: All Trades (row 30 = Trade #29) ---
$wsTrades = $wb.Worksheets.Item("All Trades")
$wsTrades.Range("G30").Value = 0.53
$wsTrades.Range("H30").Value = "CLOSED"
$wsTrades.Range("K30").Value = 100.43
$wsTrades.Range("P30").Value = "early_exit"
$wsTrades.Range("Q30").Value = 5.08

# --- Sheet: MarketMaking (row 30 = Trade #29) ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G30").Value = 0.53
$wsMM.Range("H30").Value = "CLOSED"
$wsMM.Range("K30").Value = 100.43
$wsMM.Range("P30").Value = "early_exit"
$wsMM.Range("Q30").Value = 5.08
